$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 976.7143
$ws.Range("I28").Value = 669
$ws.Range("J28").Value = 1099.8
$ws.Range("K28").Value = 669
$ws.Range("L28").Value = 1099.8
$ws.Range("M28").Value = -184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2110.7
$ws.Range("I38").Value = 3319.6
$ws.Range("J38").Value = 901.8
$ws.Range("K38").Value = 9958.799999999999
$ws.Range("L38").Value = 2705.4
$ws.Range("M38").Value = -9586.799999999999
$ws.Range("N38").Value = -3449.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 52.333332
$ws.Range("I39").Value = 39
$ws.Range("J39").Value = 99
$ws.Range("K39").Value = 117
$ws.Range("L39").Value = 297
$ws.Range("M39").Value = 179
$ws.Range("N39").Value = -889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2000
$ws.Range("I74").Value = 2000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1064

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 2000
$ws.Range("I77").Value = 2000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 10000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -5320

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2967
$ws.Range("I94").Value = 2967
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2967
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2516

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 36498.75
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 36498.75
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 36498.75
$ws.Range("N95").Value = -41990.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1882.5
$ws.Range("I98").Value = 1882.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1882.5
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -384.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1882.5
$ws.Range("I122").Value = 1882.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5647.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3197.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -6540
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2352.125
$ws.Range("I137").Value = 2302.8333
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 6908.499899999999
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -4358.499899999999
$ws.Range("N137").Value = -12600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 716
$ws.Range("I2").Value = 620
$ws.Range("J2").Value = 956
$ws.Range("K2").Value = 620
$ws.Range("L2").Value = 956
$ws.Range("M2").Value = -507
$ws.Range("N2").Value = -1182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1704.4
$ws.Range("I45").Value = 1341
$ws.Range("J45").Value = 2249.5
$ws.Range("K45").Value = 1341
$ws.Range("L45").Value = 2249.5
$ws.Range("M45").Value = -964
$ws.Range("N45").Value = -3003.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2999.6
$ws.Range("I61").Value = 2999.3333
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2999.3333
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2787.3333
$ws.Range("N61").Value = -3424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1437.1818
$ws.Range("I102").Value = 1144.2222
$ws.Range("J102").Value = 2755.5
$ws.Range("K102").Value = 1144.2222
$ws.Range("L102").Value = 2755.5
$ws.Range("M102").Value = 477.7778000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 716
$ws.Range("I116").Value = 620
$ws.Range("J116").Value = 956
$ws.Range("K116").Value = 620
$ws.Range("L116").Value = 956
$ws.Range("M116").Value = 1674
$ws.Range("N116").Value = -5544

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2999.6
$ws.Range("I136").Value = 2999.3333
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 8997.999899999999
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -6447.999899999999
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 716
$ws.Range("I3").Value = 620
$ws.Range("J3").Value = 956
$ws.Range("K3").Value = 620
$ws.Range("L3").Value = 956
$ws.Range("M3").Value = -506
$ws.Range("N3").Value = -1184

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 642
$ws.Range("I94").Value = 642
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 642
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -191
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 20000
$ws.Range("I96").Value = 20000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 20000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -17254

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2213.3333
$ws.Range("I99").Value = 2256
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2256
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -758
$ws.Range("N99").Value = -4996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5484.1665
$ws.Range("I58").Value = 2301.6667
$ws.Range("J58").Value = 8666.666999999999
$ws.Range("K58").Value = 2301.6667
$ws.Range("L58").Value = 8666.666999999999
$ws.Range("M58").Value = -2098.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5484.1665
$ws.Range("I136").Value = 2301.6667
$ws.Range("J136").Value = 8666.666999999999
$ws.Range("K136").Value = 6905.000100000001
$ws.Range("L136").Value = 26000.001
$ws.Range("M136").Value = -4355.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 114.14286
$ws.Range("I4").Value = 114.14286
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 342.42858
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -230.42858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1413.75
$ws.Range("I107").Value = 203
$ws.Range("J107").Value = 1817.3334
$ws.Range("K107").Value = 609
$ws.Range("L107").Value = 5452.0002
$ws.Range("M107").Value = 1311
$ws.Range("N107").Value = -9292.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 10357.5
$ws.Range("I129").Value = 366.66666
$ws.Range("J129").Value = 40330
$ws.Range("K129").Value = 1099.99998
$ws.Range("L129").Value = 120990
$ws.Range("M129").Value = 3900.00002
$ws.Range("N129").Value = -130990

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 833934.9
$ws.Range("I11").Value = 1300599.8
$ws.Range("J11").Value = 600602.4
$ws.Range("K11").Value = 1300599.8
$ws.Range("L11").Value = 600602.4
$ws.Range("M11").Value = -1300460.8
$ws.Range("N11").Value = -600880.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 425
$ws.Range("I29").Value = 425
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 425
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -135

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 5000
$ws.Range("I93").Value = 5000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 5000
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -3128

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 200
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 200
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -30
$ws.Range("N12").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 385
$ws.Range("I132").Value = 202.5
$ws.Range("J132").Value = 750
$ws.Range("K132").Value = 607.5
$ws.Range("L132").Value = 2250
$ws.Range("M132").Value = 1922.5
$ws.Range("N132").Value = -7310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2997.5
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 2996
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 2996
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -4244

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 2997.5
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 2996
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 14980
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -21220

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3550
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1970
